$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 46.0411815
$ws.Range("H2").Value = 92.082363
$ws.Range("I2").Value = 0.1617813835183017
$ws.Range("J2").Value = 0.1209708360422164
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.032708
$ws.Range("N2").Value = 0.065416
$ws.Range("Q2").Value = 1.505914964502
$ws.Range("R2").Value = 6.023659858008
$ws.Range("S2").Value = 0.1617813835183017
$ws.Range("T2").Value = 0.1209708360422164

# Row 3
$ws.Range("I3").Value = 0.6413247790408747
$ws.Range("J3").Value = 0.7193188085796143
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.032708
$ws.Range("N3").Value = 0.065416
$ws.Range("Q3").Value = 5.969664499465334
$ws.Range("R3").Value = 35.817986996792
$ws.Range("S3").Value = 0.6413247790408747
$ws.Range("T3").Value = 0.7193188085796143

# Row 4
$ws.Range("G4").Value = 2.911854333333334
$ws.Range("H4").Value = 8.735563000000001
$ws.Range("I4").Value = 0.01023179265394023
$ws.Range("J4").Value = 0.01147612121345596
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.032708
$ws.Range("N4").Value = 0.065416
$ws.Range("Q4").Value = 0.09524093153466669
$ws.Range("R4").Value = 0.571445589208
$ws.Range("S4").Value = 0.01023179265394023
$ws.Range("T4").Value = 0.01147612121345596

# Row 5
$ws.Range("G5").Value = 46.5307045
$ws.Range("H5").Value = 93.061409
$ws.Range("I5").Value = 0.163501489424012
$ws.Range("J5").Value = 0.1222570325437526
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.032708
$ws.Range("N5").Value = 0.065416
$ws.Range("Q5").Value = 1.521926282786
$ws.Range("R5").Value = 6.087705131144
$ws.Range("S5").Value = 0.163501489424012
$ws.Range("T5").Value = 0.1222570325437526

# Row 6
$ws.Range("G6").Value = 1.216048333333333
$ws.Range("H6").Value = 3.648145
$ws.Range("I6").Value = 0.004273000287618414
$ws.Range("J6").Value = 0.004792656663830743
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.032708
$ws.Range("N6").Value = 0.065416
$ws.Range("Q6").Value = 0.03977450888666667
$ws.Range("R6").Value = 0.23864705332
$ws.Range("S6").Value = 0.004273000287618414
$ws.Range("T6").Value = 0.004792656663830743

# Row 7
$ws.Range("G7").Value = 5.375188000000001
$ws.Range("H7").Value = 16.125564
$ws.Range("I7").Value = 0.01888755507525308
$ws.Range("J7").Value = 0.02118454495713003
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.032708
$ws.Range("N7").Value = 0.065416
$ws.Range("Q7").Value = 0.175811649104
$ws.Range("R7").Value = 1.054869894624
$ws.Range("S7").Value = 0.01888755507525308
$ws.Range("T7").Value = 0.02118454495713003
